$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 382, shifting rows 382:403 down to 383:404.
$ws.Rows.Item(382).Insert()

# Populate the newly inserted row 382 with the new weekly record.
$ws.Cells.Item(382, 1).Value = 5
$ws.Cells.Item(382, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(382, 3).Value = "Maule"
$ws.Cells.Item(382, 4).Value = 45267
$ws.Cells.Item(382, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(382, 5).Value = 7
$ws.Cells.Item(382, 6).Value = 100112017
$ws.Cells.Item(382, 7).Value = "Apio"
$ws.Cells.Item(382, 8).Value = "Americana (o)"
$ws.Cells.Item(382, 9).Value = "Primera"
$ws.Cells.Item(382, 10).Value = 300
$ws.Cells.Item(382, 11).Value = 15000
$ws.Cells.Item(382, 12).Value = 15000
$ws.Cells.Item(382, 13).Value = 15000
$ws.Cells.Item(382, 14).Value = "`$/docena de matas"
$ws.Cells.Item(382, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(382, 16).Value = 2500
$ws.Cells.Item(382, 17).Value = 6
$ws.Cells.Item(382, 18).Value = "Hortaliza"
